$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 161, shifting existing rows 161:252 down to 162:253
$ws.Rows("161").Insert()

# Populate the newly inserted row 161 with the new weekly data point
$ws.Range("A161").Value = 3
$ws.Range("B161").Value = "Femacal de La Calera"
$ws.Range("C161").Value = "Coquimbo"
$ws.Range("D161").Value = 45001
$ws.Range("E161").Value = 5
$ws.Range("F161").Value = 100112052
$ws.Range("G161").Value = "Albahaca"
$ws.Range("H161").Value = "Sin especificar"
$ws.Range("I161").Value = "Primera"
$ws.Range("J161").Value = 65
$ws.Range("K161").Value = 5000
$ws.Range("L161").Value = 5000
$ws.Range("M161").Value = 5000
$ws.Range("N161").Value = "`$/docena de matas"
$ws.Range("O161").Value = "Provincia de Quillota"
$ws.Range("P161").Value = 833
$ws.Range("Q161").Value = 6
$ws.Range("R161").Value = "Hortaliza"
